$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 606, shifting the existing rows 606..647
# down to 607..648 (dimension grows from D647 to D648).
$ws.Rows.Item(606).Insert()

# Populate the newly inserted row 606 with the new data point:
# 2026/01/12 (Monday), time slot 4, value 18.
# Force text formatting before assignment so the date-like string isn't
# auto-converted into a date serial number, then clear the formatting
# again so the cell ends up without an explicit style (matching the
# plain, unstyled cells used throughout the rest of the column).
$ws.Cells.Item(606, 1).NumberFormat = "@"
$ws.Cells.Item(606, 1).Value = "2026/01/12"
$ws.Cells.Item(606, 1).ClearFormats()

$ws.Cells.Item(606, 2).NumberFormat = "@"
$ws.Cells.Item(606, 2).Value = "月"
$ws.Cells.Item(606, 2).ClearFormats()

$ws.Cells.Item(606, 3).Value = 4
$ws.Cells.Item(606, 4).Value = 18
